$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append (device master records for two new Mac-Addresses sets)
$ids = 3000166, 3000167, 3000168, 3000169, 3000170, 3000171, 3000172, 3000173, 3000174, 3000175
$names = "Finger Print Scanner 30", "IRIS Scanner 30", "Web Camera 30", "Document Scanner 30", "Printer 30", "Finger Print Scanner 31", "IRIS Scanner 31", "Web Camera 31", "Document Scanner 31", "Printer 31"
$macs = "D6-15-AC-80-6B-86", "6D-58-E2-DF-74-34", "E2-A8-56-86-15-30", "72-E8-B9-FD-63-65", "D3-F3-A4-50-AD-12", "06-16-D0-0B-A6-E4", "21-78-45-AC-E9-20", "3C-E8-87-99-DB-FA", "BF-55-53-98-40-08", "5A-43-36-46-22-EB"
$serials = "BS563Q2230814", "BS563Q2230815", "BS563Q2230816", "BS563Q2230817", "BS563Q2230818", "BS563Q2230819", "BS563Q2230820", "BS563Q2230821", "BS563Q2230822", "BS563Q2230823"
$dspecs = 165, 327, 736, 801, 920, 165, 327, 736, 801, 920

$startRow = 147

for ($i = 0; $i -lt $ids.Count; $i++) {
    $r = $startRow + $i

    $ws.Range("A$r").Value = $ids[$i]
    $ws.Range("B$r").Value = $names[$i]
    $ws.Range("C$r").Value = $macs[$i]
    $ws.Range("D$r").Value = $serials[$i]
    $ws.Range("F$r").Value = $dspecs[$i]
    $ws.Range("G$r").Value = "eng"
    $ws.Range("H$r").Value = $true
    $ws.Range("H$r").HorizontalAlignment = -4131
    $ws.Range("I$r").Value = "superadmin"
    $ws.Range("J$r").Value = "now()"
    $ws.Range("K$r").Value = "now()"
}

# Restore the selection/view state seen in the edited workbook
$ws.Range("D145").Select() | Out-Null
